$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C corrections (existing values tweaked) ---
$ws.Cells.Item(4, 3).Value = 0.333337
$ws.Cells.Item(5, 3).Value = 0.333337
$ws.Cells.Item(8, 3).Value = 0.00001

# --- Column D: "Expected results" values, using 0.00% number format ---
$ws.Cells.Item(3, 4).Value = 0.666663
$ws.Cells.Item(4, 4).Value = 0.333333
$ws.Cells.Item(5, 4).Value = 0.333333
$ws.Cells.Item(6, 4).Value = 0.666667
$ws.Cells.Item(7, 4).Value = 0.666667
$ws.Cells.Item(8, 4).Value = 0
$ws.Range("D3:D8").NumberFormat = "0.00%"

# --- Column E: comments for each predictor ---
$ws.Cells.Item(3, 5).Value = "Main loop is always taken and it’s 1/3 of all mispredictions => 33.33% of 66.67% is main loop mispredicts. Tested branch is for() with two stages: 1 and then 0. The first is taken and the second is not taken => 1/3 is predicted and 1/3 is mispredicted. 1/3 + 1/3 = 66.67% "
$ws.Cells.Item(4, 5).Value = "Main loop is always taken and it’s 1/3 of all predictions. Tested branch is for() with two stages: 1 and then 0. The first is taken and the second is not taken => 1/3 is mispredicted and 1/3 is predicted. 1/3 = 33.37% "
$ws.Cells.Item(5, 5).Value = "Main loop is unconditional backward jump => 1/3 of all predictions is predicted. Tested branch is for() with backward jump and it consists from two stages: 1 and then 0. The first is taken and the second is not taken => 1/3 is predicted and 1/3 is mispredicted. 1/3 = 33.33%"
$ws.Cells.Item(6, 5).Value = "Main loop is always taken and predictor knows it after the first time => 1/3 predictions is true. Tested branch is jump from for and has two stages: t0 = 1 – jump, t0 = 0 – not jump. => NT becomes T after 1, but T becomes NT after 0. It changes after each prediction  and it always fail => 2/3 of all predictions is fail. 66.67%"
$ws.Cells.Item(7, 5).Value = "Main loop is always taken and predictor knows it after the first time => 1/3 predictions is true. Tested branch is jump from for and has two stages: t0 = 1 – jump, t0 = 0 – not jump. => WEAKLY NT becomes WEAKLY T after 1, but WEAKLY T becomes WEAKLY NT after 0. It changes after each prediction  and it always fail => 2/3 of all predictions is fail. 66.67%(it isn’t better then one bit, because we can’t use strongly NT and T)"
$ws.Cells.Item(8, 5).Value = "Main loop is always taken and predictor knows it after the first two times => 1/3 of all predictions is true(0 misses). Tested branch is for() with two stages: t0 = 1 and then t0 = 0 (The first is taken and the second is not taken) => we have history 101010101010… with two patterns 01 and 10. => after 01 is always 0 and after 10 is always 1. Predictor can predict all stages(except first two) => it is always true. => 0% mispredicts"

# --- Selection moves to B5 ---
$ws.Range("B5").Select()

